$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before the current row 120, shifting the
# existing rows 120-122 down to 122-124 (cell formatting/styles carry
# down with the shifted rows, matching the target dimension A1:R124).
$ws.Rows("120:121").Insert()

# Populate the newly inserted row 120 with its final values.
$ws.Cells.Item(120, 1).Value = 6
$ws.Cells.Item(120, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(120, 3).Value = "Metropolitana"
$ws.Cells.Item(120, 4).Value = 44448
$ws.Cells.Item(120, 5).Value = 13
$ws.Cells.Item(120, 6).Value = 100112022
$ws.Cells.Item(120, 7).Value = "Arveja Verde"
$ws.Cells.Item(120, 8).Value = "Perfection"
$ws.Cells.Item(120, 9).Value = "Primera"
$ws.Cells.Item(120, 10).Value = 350
$ws.Cells.Item(120, 11).Value = 34000
$ws.Cells.Item(120, 12).Value = 35000
$ws.Cells.Item(120, 13).Value = 34571
$ws.Cells.Item(120, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(120, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(120, 16).Value = 1383
$ws.Cells.Item(120, 17).Value = 25
$ws.Cells.Item(120, 18).Value = "Hortaliza"

# Populate the newly inserted row 121 with its final values.
$ws.Cells.Item(121, 1).Value = 6
$ws.Cells.Item(121, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(121, 3).Value = "Metropolitana"
$ws.Cells.Item(121, 4).Value = 44448
$ws.Cells.Item(121, 5).Value = 13
$ws.Cells.Item(121, 6).Value = 100112022
$ws.Cells.Item(121, 7).Value = "Arveja Verde"
$ws.Cells.Item(121, 8).Value = "Sin especificar"
$ws.Cells.Item(121, 9).Value = "Primera"
$ws.Cells.Item(121, 10).Value = 430
$ws.Cells.Item(121, 11).Value = 30000
$ws.Cells.Item(121, 12).Value = 32000
$ws.Cells.Item(121, 13).Value = 31209
$ws.Cells.Item(121, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(121, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(121, 16).Value = 1248
$ws.Cells.Item(121, 17).Value = 25
$ws.Cells.Item(121, 18).Value = "Hortaliza"
